$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update vm_pu values for the 380 kV case (Case_4_22)
# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.06719484907255
$ws.Range("D2").Value = 1.069666558488837
$ws.Range("E2").Value = 1.079632953714957
$ws.Range("F2").Value = 1.084202344188526
$ws.Range("I2").Value = 1.047412803875399
$ws.Range("J2").Value = 1.072140988806746
$ws.Range("K2").Value = 1.072368455464166
$ws.Range("L2").Value = 1.082308494517143
$ws.Range("M2").Value = 1.086865976127654
$ws.Range("N2").Value = 1.027552197144529

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.068595153063681
$ws.Range("D3").Value = 1.070771339435266
$ws.Range("E3").Value = 1.080964263742628
$ws.Range("F3").Value = 1.085480160068898
$ws.Range("I3").Value = 1.047735809689218
$ws.Range("J3").Value = 1.07319544925821
$ws.Range("K3").Value = 1.073288735268127
$ws.Range("L3").Value = 1.083456664211223
$ws.Range("M3").Value = 1.087961647803132
$ws.Range("N3").Value = 1.02791720160511

# Row 4
$ws.Range("B4").Value = 1.019999999999999
$ws.Range("C4").Value = 1.069500669276834
$ws.Range("D4").Value = 1.071485482384905
$ws.Range("E4").Value = 1.081825421620733
$ws.Range("F4").Value = 1.086306577817918
$ws.Range("I4").Value = 1.047943159619019
$ws.Range("J4").Value = 1.073876699884925
$ws.Range("K4").Value = 1.073882900476671
$ws.Range("L4").Value = 1.084198767171872
$ws.Range("M4").Value = 1.088669643243165
$ws.Range("N4").Value = 1.028152709191308

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.069881215244747
$ws.Range("D5").Value = 1.07178553735695
$ws.Range("E5").Value = 1.082187386745757
$ws.Range("F5").Value = 1.086653907794574
$ws.Range("I5").Value = 1.04802993398792
$ws.Range("J5").Value = 1.074162847909408
$ws.Range("K5").Value = 1.074132374675541
$ws.Range("L5").Value = 1.084510549406263
$ws.Range("M5").Value = 1.088967053475886
$ws.Range("N5").Value = 1.028251555949186

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.069945102912729
$ws.Range("D6").Value = 1.071835907972554
$ws.Range("E6").Value = 1.082248158534231
$ws.Range("F6").Value = 1.086712220486772
$ws.Range("I6").Value = 1.048044480620845
$ws.Range("J6").Value = 1.074210878864615
$ws.Range("K6").Value = 1.074174244218498
$ws.Range("L6").Value = 1.084562887504747
$ws.Range("M6").Value = 1.089016976499014
$ws.Range("N6").Value = 1.028268143361204

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.069505754667999
$ws.Range("D7").Value = 1.071489492401665
$ws.Range("E7").Value = 1.081830258470071
$ws.Range("F7").Value = 1.086311219231473
$ws.Range("I7").Value = 1.047944320655338
$ws.Range("J7").Value = 1.073880524386385
$ws.Range("K7").Value = 1.073886235189302
$ws.Range("L7").Value = 1.084202933994844
$ws.Range("M7").Value = 1.088673618157468
$ws.Range("N7").Value = 1.028154030616585

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.067668209882686
$ws.Range("D8").Value = 1.070040075930721
$ws.Range("E8").Value = 1.080082938051498
$ws.Range("F8").Value = 1.084634275819556
$ws.Range("I8").Value = 1.047522308787432
$ws.Range("J8").Value = 1.072497569083788
$ws.Range("K8").Value = 1.072679742623571
$ws.Range("L8").Value = 1.082696699588505
$ws.Range("M8").Value = 1.08723646765083
$ws.Range("N8").Value = 1.027675692353539

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.064425627166142
$ws.Range("D9").Value = 1.067480353150963
$ws.Range("E9").Value = 1.077001568732715
$ws.Range("F9").Value = 1.081675968133733
$ws.Range("I9").Value = 1.046765938095452
$ws.Range("J9").Value = 1.070052403299627
$ws.Range("K9").Value = 1.070543545601639
$ws.Range("L9").Value = 1.080035950034536
$ws.Range("M9").Value = 1.084696406904557
$ws.Range("N9").Value = 1.026827591409824

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.062260539297211
$ws.Range("D10").Value = 1.065769892795016
$ws.Range("E10").Value = 1.074945510213756
$ws.Range("F10").Value = 1.079701327709284
$ws.Range("I10").Value = 1.046253065874998
$ws.Range("J10").Value = 1.068416573161213
$ws.Range("K10").Value = 1.069112393542026
$ws.Range("L10").Value = 1.078257498099773
$ws.Range("M10").Value = 1.082997724179572
$ws.Range("N10").Value = 1.0262586318833

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.061322171719622
$ws.Range("D11").Value = 1.065028262552892
$ws.Range("E11").Value = 1.074054734565939
$ws.Range("F11").Value = 1.078845662937313
$ws.Range("I11").Value = 1.046028925961745
$ws.Range("J11").Value = 1.067706840766823
$ws.Range("K11").Value = 1.068490987259067
$ws.Range("L11").Value = 1.0774862680407
$ws.Range("M11").Value = 1.082260874320024
$ws.Range("N11").Value = 1.026011408894987

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.060973484033163
$ws.Range("D12").Value = 1.064752635979335
$ws.Range("E12").Value = 1.073723783025517
$ws.Range("F12").Value = 1.078527731865998
$ws.Range("I12").Value = 1.045945359158047
$ws.Range("J12").Value = 1.067442999331639
$ws.Range("K12").Value = 1.068259909461559
$ws.Range("L12").Value = 1.07719962230507
$ws.Range("M12").Value = 1.081986975071418
$ws.Range("N12").Value = 1.025919448927037

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.061048284975991
$ws.Range("D13").Value = 1.064811765751793
$ws.Range("E13").Value = 1.073794776845628
$ws.Range("F13").Value = 1.078595933708887
$ws.Range("I13").Value = 1.04596329862787
$ws.Range("J13").Value = 1.067499604015878
$ws.Range("K13").Value = 1.068309488243691
$ws.Range("L13").Value = 1.077261116841966
$ws.Range("M13").Value = 1.08204573651885
$ws.Range("N13").Value = 1.025939180569525

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.061293351896439
$ws.Range("D14").Value = 1.065005482296481
$ws.Range("E14").Value = 1.074027379620912
$ws.Range("F14").Value = 1.078819384700471
$ws.Range("I14").Value = 1.046022024661066
$ws.Range("J14").Value = 1.067685035964106
$ws.Range("K14").Value = 1.068471891615552
$ws.Range("L14").Value = 1.077462577439507
$ws.Range("M14").Value = 1.08223823783252
$ws.Range("N14").Value = 1.026003810121111

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.061444327504957
$ws.Range("D15").Value = 1.065124817198881
$ws.Range("E15").Value = 1.074170683426698
$ws.Range("F15").Value = 1.078957046980773
$ws.Range("I15").Value = 1.046058166423156
$ws.Range("J15").Value = 1.067799258073837
$ws.Range("K15").Value = 1.068571919165276
$ws.Range("L15").Value = 1.07758668044824
$ws.Range("M15").Value = 1.082356817577848
$ws.Range("N15").Value = 1.026043613214007

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.062322796863851
$ws.Range("D16").Value = 1.065819091211661
$ws.Range("E16").Value = 1.075004617389662
$ws.Range("F16").Value = 1.079758101691369
$ws.Range("I16").Value = 1.046267897726691
$ws.Range("J16").Value = 1.068463645790598
$ws.Range("K16").Value = 1.069153597969765
$ws.Range("L16").Value = 1.07830865755887
$ws.Range("M16").Value = 1.083046598579844
$ws.Range("N16").Value = 1.026275021044713

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.062873600659675
$ws.Range("D17").Value = 1.066254323612451
$ws.Range("E17").Value = 1.075527588200291
$ws.Range("F17").Value = 1.080260409966383
$ws.Range("I17").Value = 1.046398903403631
$ws.Range("J17").Value = 1.068880019231119
$ws.Range("K17").Value = 1.069518010162593
$ws.Range("L17").Value = 1.07876122439767
$ws.Range("M17").Value = 1.083478926882697
$ws.Range("N17").Value = 1.026419946055192

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.063194791598941
$ws.Range("D18").Value = 1.066508092055993
$ws.Range("E18").Value = 1.075832581367434
$ws.Range("F18").Value = 1.080553337140155
$ws.Range("I18").Value = 1.046475117853996
$ws.Range("J18").Value = 1.069122747304499
$ws.Range("K18").Value = 1.069730401052205
$ws.Range("L18").Value = 1.079025088422302
$ws.Range("M18").Value = 1.083730970619199
$ws.Range("N18").Value = 1.026504395487666

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.063304295393488
$ws.Range("D19").Value = 1.066594604483032
$ws.Range("E19").Value = 1.075936568339246
$ws.Range("F19").Value = 1.080653207555217
$ws.Range("I19").Value = 1.046501071307787
$ws.Range("J19").Value = 1.069205488434773
$ws.Range("K19").Value = 1.069802793038821
$ws.Range("L19").Value = 1.079115040603501
$ws.Range("M19").Value = 1.083816889752328
$ws.Range("N19").Value = 1.026533176540096

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.062814513312805
$ws.Range("D20").Value = 1.066207637173033
$ws.Range("E20").Value = 1.075471483261531
$ws.Range("F20").Value = 1.080206523356098
$ws.Range("I20").Value = 1.046384868322276
$ws.Range("J20").Value = 1.068835360346717
$ws.Range("K20").Value = 1.069478929230493
$ws.Range("L20").Value = 1.078712679724301
$ws.Range("M20").Value = 1.083432555213324
$ws.Range("N20").Value = 1.026404405554112

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.061221189599776
$ws.Range("D21").Value = 1.064948441827911
$ws.Range("E21").Value = 1.073958886084058
$ws.Range("F21").Value = 1.078753586720532
$ws.Range("I21").Value = 1.046004739911091
$ws.Range("J21").Value = 1.067630436849067
$ws.Range("K21").Value = 1.068424075069097
$ws.Range("L21").Value = 1.077403257220615
$ws.Range("M21").Value = 1.082181556555804
$ws.Range("N21").Value = 1.025984781936571

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.060218613091144
$ws.Range("D22").Value = 1.06415585414724
$ws.Range("E22").Value = 1.073007403876401
$ws.Range("F22").Value = 1.07783949253828
$ws.Range("I22").Value = 1.045763936568412
$ws.Range("J22").Value = 1.066871606505236
$ws.Range("K22").Value = 1.067759341448693
$ws.Range("L22").Value = 1.076578947854964
$ws.Range("M22").Value = 1.081393842716292
$ws.Range("N22").Value = 1.025720193205284

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.060750174635397
$ws.Range("D23").Value = 1.064576104577654
$ws.Range("E23").Value = 1.073511846938622
$ws.Range("F23").Value = 1.078324127047281
$ws.Range("I23").Value = 1.045891762188827
$ws.Range("J23").Value = 1.067273996318479
$ws.Range("K23").Value = 1.068111873051161
$ws.Range("L23").Value = 1.077016028193539
$ws.Range("M23").Value = 1.081811536059194
$ws.Range("N23").Value = 1.025860528629084

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.062841212599466
$ws.Range("D24").Value = 1.066228733059603
$ws.Range("E24").Value = 1.075496834812941
$ws.Range("F24").Value = 1.080230872583289
$ws.Range("I24").Value = 1.046391210786064
$ws.Range("J24").Value = 1.06885554019278
$ws.Range("K24").Value = 1.069496588729866
$ws.Range("L24").Value = 1.078734615314031
$ws.Range("M24").Value = 1.083453508963973
$ws.Range("N24").Value = 1.026411427893648

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.065264485417703
$ws.Range("D25").Value = 1.068142791623984
$ws.Range("E25").Value = 1.07779847976161
$ws.Range("F25").Value = 1.082441175823165
$ws.Range("I25").Value = 1.04696299349995
$ws.Range("J25").Value = 1.070685531275965
$ws.Range("K25").Value = 1.071097029343156
$ws.Range("L25").Value = 1.080724616948745
$ws.Range("M25").Value = 1.085353995452613
$ws.Range("N25").Value = 1.027047468876305

